$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (pushes existing rows 4-9 down to 5-10)
$ws.Rows.Item(3).Insert()

# Header row (now row 5) - add Clue columns
$ws.Range("D5").Value = "Clue 1"
$ws.Range("E5").Value = "Clue 2"
$ws.Range("F5").Value = "Clue 3"

# New note cell A3
$ws.Range("A3").Value = "Minimum of 3000 PA career leaderboards for rate statistics."
$ws.Range("A3").Font.Name = "Calibri"

# Data rows (now rows 6-8) - add Clue data
$ws.Range("D6").Value = "12 Batting Titles"
$ws.Range("E6").Value = "1 Triple Crown and 1 MVP win"
$ws.Range("F6").Value = "Nicknamed the Georgia Peach"

$ws.Range("D7").Value = "7 Batting Titles"
$ws.Range("E7").Value = "2 Triple Crown and 7 MVP wins"
$ws.Range("F7").Value = "Nicknamed Rajah"

$ws.Range("D8").Value = "13 year career"
$ws.Range("E8").Value = "Played for PHA, CLE, and CWS"
$ws.Range("F8").Value = "3 time AL Triples leader"

# Update formulas for C6, C7 (the AVG column), replacing static values with formulas
$ws.Range("C6").Formula = "=4189/11439"
$ws.Range("C7").Formula = "=2930/8173"

# Copy style from C5 (header style) to D5:F5
$ws.Range("C5").Copy()
$ws.Range("D5:F5").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("D9").Select()
